$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.514.79"
$ws.Range("E2").Value = "  +2.72%  "
$ws.Range("D3").Value = "2.339.48"
$ws.Range("E3").Value = "  +1.85%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.77"
$ws.Range("E5").Value = "  +6.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.36"
$ws.Range("E6").Value = "  +3.20%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.570"
$ws.Range("E8").Value = "  +7.87%  "
$ws.Range("D9").Value = "2.338.04"
$ws.Range("E9").Value = "  +1.73%  "
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("E12").Value = "  +3.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.358"
$ws.Range("E13").Value = "  +7.41%  "
$ws.Range("D14").Value = "2.758.10"
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.54"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "57.513.84"
$ws.Range("E16").Value = "  +2.68%  "
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "2.347.21"
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("E19").Value = "  +3.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "333.43"
$ws.Range("E20").Value = "  +1.87%  "
$ws.Range("E21").Value = "  +2.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.69"
$ws.Range("E22").Value = "  +0.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.57"
$ws.Range("E24").Value = "  +0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "62.08"
$ws.Range("E25").Value = "  +2.64%  "
$ws.Range("E26").Value = "  +2.35%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("E28").Value = "  -1.05%  "
$ws.Range("E29").Value = "  +6.59%  "
$ws.Range("E30").Value = "  +5.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "170.19"
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("D32").Value = "0.0₃0732"
$ws.Range("E32").Value = "  +3.63%  "
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("E34").Value = "  +18.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.47"
$ws.Range("E35").Value = "  +1.56%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.16"
$ws.Range("E38").Value = "  +7.41%  "
$ws.Range("E39").Value = "  +1.67%  "
$ws.Range("E40").Value = "  +3.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.99"
$ws.Range("E41").Value = "  +1.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "147.52"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.376"
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("E44").Value = "  +1.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "283.08"
$ws.Range("E45").Value = "  +2.20%  "
$ws.Range("E46").Value = "  +2.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.11"
$ws.Range("E47").Value = "  +6.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0505"
$ws.Range("E48").Value = "  +2.76%  "
$ws.Range("E49").Value = "  +1.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.384"
$ws.Range("E50").Value = "  +1.97%  "
$ws.Range("E51").Value = "  +2.05%  "
